$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.086.52"
$ws.Range("E2").Value = "  -1.00%  "

# Row 3
$ws.Range("D3").Value = "1.666.65"
$ws.Range("E3").Value = "  -1.24%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.37%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5171"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.53%  "

# Row 7
$ws.Range("E7").Value = "  -0.46%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2643"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06198"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.99%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.09%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07491"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.47%  "

# Row 12
$ws.Range("D12").Value = "1.670.33"
$ws.Range("E12").Value = "  -1.08%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.429"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.34%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5581"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.98%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000007951"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.47%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.81%  "

# Row 17
$ws.Range("D17").Value = "26.081.18"
$ws.Range("E17").Value = "  -1.25%  "

# Row 18
$ws.Range("E18").Value = "  -0.45%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.800"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.86%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.43%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "185.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.53%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.157"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.17%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.008"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.24%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "146.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.14%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1247"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.85%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.548"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.22%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06291"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.86%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.347"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.67%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.271"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.487"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.51%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.435"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.93%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.627"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.80%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9949"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.63%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.415"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6009"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.71%  "

# Row 37
$ws.Range("E37").Value = "  -0.63%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.103"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.13%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01603"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.67%  "

# Row 40
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.093.33"
$ws.Range("E40").Value = "  -2.00%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8585"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.23%  "

# Row 42
$ws.Range("E42").Value = "  -1.08%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.46%  "

# Row 44
$ws.Range("D44").Value = "1.823.66"
$ws.Range("E44").Value = "  -0.88%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000106"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.53%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.98%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9958"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.74%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05247"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.46%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.920"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.50%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4269"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.82%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.892"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.41%  "
